# Regenerate the "K" column (column G) on Sheet1 using updated values
# (commit: "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new K value (column G), per the canonical diff.
$kValues = @{
    2  = 4
    3  = 3
    4  = 4
    5  = 3
    6  = 5
    7  = 0
    8  = 4
    9  = 1
    10 = 1
    11 = 2
    12 = 3
    13 = 2
    14 = 3
    15 = 4
    16 = 1
    17 = 1
    18 = 0
    19 = 0
    20 = 1
    21 = 0
    22 = 4
    23 = 1
    24 = 0
    25 = 5
    26 = 2
    27 = 4
    28 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
